$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.438.27'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '1.562.92'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = "'208.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.13%  '

$ws.Range("E6").Value = '  -0.46%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = "'21.87"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.90%  '

$ws.Range("E9").Value = '  -0.94%  '

$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.30%  '

$ws.Range("D12").Value = '1.785.03'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").Value = '1.563.94'
$ws.Range("E13").Value = '  -0.75%  '

$ws.Range("E14").Value = '  -0.39%  '

$ws.Range("D15").Value = "'0.515"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("E16").Value = '  +1.13%  '

$ws.Range("D17").Value = '27.418.60'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").Value = "'213.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.56%  '

$ws.Range("D19").Value = '0.0₃0687'
$ws.Range("E19").Value = '  -0.38%  '

$ws.Range("E20").Value = '  -0.79%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = "'4.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("E23").Value = '  +1.37%  '

$ws.Range("E24").Value = '  +1.77%  '

$ws.Range("D25").Value = "'152.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.42%  '

$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  +1.09%  '

$ws.Range("D28").Value = "'15.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("E29").Value = '  -1.34%  '

$ws.Range("E30").Value = '  +0.33%  '

$ws.Range("E31").Value = '  +1.59%  '

$ws.Range("E32").Value = '  -1.07%  '

$ws.Range("D33").Value = '1.361.16'
$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("E35").Value = '  +1.62%  '

$ws.Range("D36").Value = "'0.973"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.94%  '

$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("E38").Value = '  +2.06%  '

$ws.Range("D39").Value = "'0.532"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.56%  '

$ws.Range("D40").Value = "'0.820"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.02%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("E43").Value = '  +1.81%  '

$ws.Range("D44").Value = "'64.08"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.35%  '

$ws.Range("D45").Value = "'5.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.31%  '

$ws.Range("E46").Value = '  -1.16%  '

$ws.Range("D47").Value = '1.698.80'
$ws.Range("E47").Value = '  -0.65%  '

$ws.Range("D48").Value = "'85.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.02%  '

$ws.Range("D49").Value = '0.0₇0988'
$ws.Range("E49").Value = '  -1.04%  '

$ws.Range("D50").Value = "'0.0954"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.16%  '

$ws.Range("D51").Value = "'0.0493"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.62%  '
